$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7903876304626465
$ws.Range("B1").Value = 1.234026908874512
$ws.Range("C1").Value = 4.472187995910645
$ws.Range("D1").Value = 4.072829246520996
$ws.Range("E1").Value = 1.09970486164093
